$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 7 and 8, pushing the existing "Id" row (old row 7) down to row 9.
$ws.Rows("7:8").Insert()

# Row 7: Icon / string / Friend relation / 图标 description
$ws.Range("A7").Value = "Icon"
$ws.Range("B7").Value = "string"
$ws.Range("C7").Value = $false
$ws.Range("D7").Value = $false
$ws.Range("E7").Value = $false
$ws.Range("F7").Value = $false
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = "Friend"
$ws.Range("J7").Value = "图标"

# Row 8: ShowName / string / Friend relation / 名字 description
$ws.Range("A8").Value = "ShowName"
$ws.Range("B8").Value = "string"
$ws.Range("C8").Value = $false
$ws.Range("D8").Value = $false
$ws.Range("E8").Value = $false
$ws.Range("F8").Value = $false
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = "Friend"
$ws.Range("J8").Value = "名字"

# Restore the active selection cell as it appears in the edited workbook.
$ws.Range("C13").Select()
